$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-change (E) columns for each coin row.
# D-column values are plain inline-string text in the workbook (not numbers),
# so we momentarily force a Text number format while writing the value and
# then restore the original cell style, avoiding any numeric auto-coercion
# or stray style-table changes.

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "71.506.62"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  +3.10%  "
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.000.23"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  -0.17%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "528.32"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  +3.90%  "
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "148.86"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  +1.37%  "
$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.624"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +0.19%  "
$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.735"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  +0.42%  "
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.176"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  +1.49%  "
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0000345"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  +0.03%  "
$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "44.06"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  +1.59%  "
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "10.65"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  +1.78%  "
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.647.28"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  +1.45%  "
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.011.42"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  +0.74%  "
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "21.36"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  +7.41%  "
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "14.29"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("E18").Value = "  -0.62%  "
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.133"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  -1.82%  "
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "71.429.62"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  +2.61%  "
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "441.08"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  +1.67%  "
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.59"
$cell.Style = $origStyle
$ws.Range("E22").Value = "  +4.68%  "
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "93.53"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  +5.57%  "
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "14.36"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  -1.36%  "
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "12.29"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  +4.43%  "
$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.11"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  +6.21%  "
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "10.86"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  -2.37%  "
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "36.94"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  -0.11%  "
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "13.64"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  +2.34%  "
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "700.86"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("E32").Value = "  -0.01%  "
$cell = $ws.Range("D33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.95"
$cell.Style = $origStyle
$ws.Range("E33").Value = "  +15.98%  "
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "67.26"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  -0.96%  "
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0912"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  +3.96%  "
$ws.Range("E36").Value = "  +0.01%  "
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "41.15"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  +1.10%  "
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.59"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  +17.92%  "
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.152"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("E40").Value = "  +0.13%  "
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0496"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("E42").Value = "  -0.14%  "
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.90"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  +1.10%  "
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.15"
$cell.Style = $origStyle
$ws.Range("E44").Value = "  -0.19%  "
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.54"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +5.47%  "
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.23"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  +8.42%  "
$ws.Range("E47").Value = "  +1.45%  "
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.000287"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  +21.89%  "
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.26"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  +6.00%  "
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.40"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  +0.83%  "
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0₆0344"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  -5.63%  "
